# Update TPM-derived NATMI metrics (columns G:T) for rows 2-10 on the
# active worksheet with newly recomputed values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(20.09599933333334, 60.287998, 0.9230842418515092, 0.9230842418515091, 3, 1, 0.1869266666666667, 0.5607800000000001, 0.05045642741357319, 0.05045642741357319, 3.756478168715556, 33.80830351844001, 0.04657553304559391, 0.04657553304559391)
    3  = @(20.09599933333334, 60.287998, 0.9230842418515092, 0.9230842418515091, 3, 1, 2.050040666666666, 6.150122, 0.5533599348721772, 0.5533599348721773, 41.19761587063955, 370.778542835756, 0.5107978359524842, 0.5107978359524842)
    4  = @(20.09599933333334, 60.287998, 0.9230842418515092, 0.9230842418515091, 3, 1, 1.467747333333334, 4.403242000000001, 0.3961836377142495, 0.3961836377142496, 29.49584943216845, 265.462644889516, 0.365710872853431, 0.365710872853431)
    5  = @(1.158105666666667, 3.474317, 0.05319611498621682, 0.05319611498621682, 3, 1, 0.1869266666666667, 0.5607800000000001, 0.05045642741357319, 0.05045642741357319, 0.2164808319177778, 1.94832748726, 0.002684085914486142, 0.002684085914486142)
    6  = @(1.158105666666667, 3.474317, 0.05319611498621682, 0.05319611498621682, 3, 1, 2.050040666666666, 6.150122, 0.5533599348721772, 0.5533599348721773, 2.374163712963778, 21.367473416674, 0.02943659872422579, 0.02943659872422579)
    7  = @(1.158105666666667, 3.474317, 0.05319611498621682, 0.05319611498621682, 3, 1, 1.467747333333334, 4.403242000000001, 0.3961836377142495, 0.3961836377142496, 1.699806503968222, 15.298258535714, 0.02107543034750489, 0.02107543034750489)
    8  = @(0.5163883333333333, 1.549165, 0.02371964316227407, 0.02371964316227407, 3, 1, 0.1869266666666667, 0.5607800000000001, 0.05045642741357319, 0.05045642741357319, 0.09652674985555557, 0.8687407487000001, 0.00119680845349314, 0.001196808453493139)
    9  = @(0.5163883333333333, 1.549165, 0.02371964316227407, 0.02371964316227407, 3, 1, 2.050040666666666, 6.150122, 0.5533599348721772, 0.5533599348721773, 1.058617083125555, 9.527553748129998, 0.01312550019546726, 0.01312550019546727)
    10 = @(0.5163883333333333, 1.549165, 0.02371964316227407, 0.02371964316227407, 3, 1, 1.467747333333334, 4.403242000000001, 0.3961836377142495, 0.3961836377142496, 0.7579275992144445, 6.82134839293, 0.009397334513313668, 0.009397334513313668)
}

# Columns G (7) through T (20)
$startCol = 7

foreach ($rowNum in $data.Keys) {
    $values = $data[$rowNum]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $ws.Cells.Item($rowNum, $startCol + $i).Value = $values[$i]
    }
}
